# Add 2022-Q4 data
# -----------------------------------------------------------------------
# 1) Insert a brand-new worksheet named "2022-Q4" right after the "总计"
#    (summary) sheet and before the existing "2022-Q2" sheet, and fill it
#    with the quarterly fund-holding detail rows.
# 2) Update the "总计" (summary) sheet so that it gains a new top data row
#    for 2022-Q4 and the existing 2021-Q4/2022-Q1/2022-Q2 rows shift down
#    to make room, with a new row appended for the (previously last)
#    2021-Q4 entry.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Create the new "2022-Q4" sheet, positioned before the current "2022-Q2"
# sheet (i.e. directly after "总计").
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($q2Sheet)
$q4Sheet.Name = "2022-Q4"

# Header row
$q4Sheet.Cells.Item(1, 2).Value = "基金代码"
$q4Sheet.Cells.Item(1, 3).Value = "基金名称"
$q4Sheet.Cells.Item(1, 4).Value = "基金规模"
$q4Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q4Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q4Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4Sheet.Cells.Item(1, 8).Value = "仓位排名"

for ($col = 2; $col -le 8; $col++) {
    $q4Sheet.Cells.Item(1, $col).Font.Bold = $true
    $q4Sheet.Cells.Item(1, $col).HorizontalAlignment = -4108
    $q4Sheet.Cells.Item(1, $col).VerticalAlignment = -4160
    $q4Sheet.Cells.Item(1, $col).Borders.LineStyle = 1
}

# Data rows
$q4Data = @(
    @(0, "005051", "上投摩根标普港股通低波红利指数A", "1.64", "93.98", "2.38", "0.0390", 8),
    @(1, "005052", "上投摩根标普港股通低波红利指数C", "1.52", "93.98", "2.38", "0.0362", 8),
    @(2, "005702", "恒生前海港股通高股息低波动指数", "0.23", "94.47", "2.31", "0.0053", 9)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = $q4Data[$i]
    $r = $i + 2
    $q4Sheet.Cells.Item($r, 1).Value = $row[0]
    $q4Sheet.Cells.Item($r, 1).Font.Bold = $true
    $q4Sheet.Cells.Item($r, 1).HorizontalAlignment = -4108
    $q4Sheet.Cells.Item($r, 1).VerticalAlignment = -4160
    $q4Sheet.Cells.Item($r, 1).Borders.LineStyle = 1

    $q4Sheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4Sheet.Cells.Item($r, 3).Value = $row[2]
    $q4Sheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4Sheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4Sheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4Sheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4Sheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# Update the "总计" (summary) sheet: row 2 becomes the new 2022-Q4 entry,
# rows 3 and 4 shift to hold what used to be in rows 2 and 3, and a new
# row 5 is appended with the former row-4 (2021-Q4) data.
# ---------------------------------------------------------------------
$summaryData = @(
    @("2022-Q4", 3, 0.08),
    @("2022-Q2", 2, 0.09),
    @("2022-Q1", 3, 0.5),
    @("2021-Q4", 3, 0.07000000000000001)
)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $row = $summaryData[$i]
    $r = $i + 2
    $summarySheet.Cells.Item($r, 1).Value = $i
    $summarySheet.Cells.Item($r, 1).Font.Bold = $true
    $summarySheet.Cells.Item($r, 1).HorizontalAlignment = -4108
    $summarySheet.Cells.Item($r, 1).VerticalAlignment = -4160
    $summarySheet.Cells.Item($r, 1).Borders.LineStyle = 1

    $summarySheet.Cells.Item($r, 2).Value = $row[0]
    $summarySheet.Cells.Item($r, 3).Value = $row[1]
    $summarySheet.Cells.Item($r, 4).Value = $row[2]
}

# Restore the originally active/selected sheet ("2021-Q4", now the last
# tab) so the tab-selection state matches the source workbook.
$wb.Worksheets.Item("2021-Q4").Activate()

